$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-14 Saturday" "2024-12-15 Sunday"

Replace-Text "642×3=" "767×8="
Replace-Text "643×4=" "917×6="
Replace-Text "683×8=" "376×4="
Replace-Text "392×5=" "826×5="
Replace-Text "917×3=" "937×6="

Replace-Text "773×8=" "280×7="
Replace-Text "298×8=" "852×6="
Replace-Text "137×6=" "352×4="
Replace-Text "354×2=" "153×3="
Replace-Text "565×2=" "693×6="

Replace-Text "841×4=" "917×6="
Replace-Text "347×6=" "649×6="
Replace-Text "899×8=" "428×6="
Replace-Text "298×4=" "172×5="
Replace-Text "369×5=" "954×8="

Replace-Text "393×8=" "877×8="
Replace-Text "792×9=" "342×7="
Replace-Text "434×4=" "972×5="
Replace-Text "481×8=" "919×3="
Replace-Text "763×6=" "596×2="

Replace-Text "508×4=" "559×6="
Replace-Text "436×8=" "494×4="
Replace-Text "142×5=" "990×6="
Replace-Text "590×7=" "164×8="
Replace-Text "290×6=" "864×6="
